# Updates to mapping; temporary fix for publication author parsing
#
# 1. Fix the "Docotr (other health professions)" typo on the degree-type
#    sheet (now spelled correctly: "Doctor (other health professions)").
# 2. Add a new "c-placement-type" choice-list sheet (Academic / Corporate /
#    Government / Other) after "degree-type".

$wb = $excel.ActiveWorkbook
$degreeType = $wb.Worksheets.Item("degree-type")

# --- 1. Fix the misspelled "Docotr (other health professions)" entry -----
$fixed = "Doctor (other health professions)"
for ($r = 1; $r -le $degreeType.UsedRange.Rows.Count; $r++) {
    $cell = $degreeType.Cells.Item($r, 1)
    if ($cell.Value() -eq "Docotr (other health professions)") {
        $cell.Value = $fixed
    }
}

# --- 2. Add the new c-placement-type sheet --------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$placementType = $wb.Worksheets.Add($null, $lastSheet)
$placementType.Name = "c-placement-type"

# Reuse the existing "Arial 10" body style already used on degree-type
# (rather than building the font up property-by-property) so the new
# sheet's formatting matches the rest of the workbook.
$degreeType.Range("A2").Copy()
$placementType.Range("A1:B5").PasteSpecial(-4122) # xlPasteFormats

# Header row - values, then bump to bold + underlined for emphasis.
$placementType.Cells.Item(1, 1).Value = "Elements"
$placementType.Cells.Item(1, 2).Value = "Lyterati"
$header = $placementType.Range("A1:B1")
$header.Font.Bold = $true
$header.Font.Underline = $true

# Data rows - both columns map to the same literal choice value.
$values = @("Academic", "Corporate", "Government", "Other")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $placementType.Cells.Item($row, 1).Value = $values[$i]
    $placementType.Cells.Item($row, 2).Value = $values[$i]
}

$placementType.Range("D3").Select() | Out-Null
$degreeType.Select() | Out-Null
$degreeType.Range("E9").Select() | Out-Null
